# Updates profit-sheet price/profit figures across the 8 crafting-job
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed
# market-board pricing pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3421.9119
$ws.Range("I137").Value = 2366.2903
$ws.Range("K137").Value = 7098.8709
$ws.Range("M137").Value = -4548.8709
$ws.Range("H138").Value = 2277.3699
$ws.Range("J138").Value = 2404.34
$ws.Range("L138").Value = 7213.02
$ws.Range("N138").Value = -17493.02

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1779
$ws.Range("I2").Value = 1779
$ws.Range("K2").Value = 1779
$ws.Range("M2").Value = -1666
$ws.Range("H61").Value = 2782.5151
$ws.Range("I61").Value = 2384.8096
$ws.Range("J61").Value = 3478.5
$ws.Range("K61").Value = 2384.8096
$ws.Range("L61").Value = 3478.5
$ws.Range("M61").Value = -2172.8096
$ws.Range("N61").Value = -3902.5
$ws.Range("H74").Value = 2875.7646
$ws.Range("I74").Value = 2726.7144
$ws.Range("J74").Value = 3571.3333
$ws.Range("K74").Value = 2726.7144
$ws.Range("L74").Value = 3571.3333
$ws.Range("M74").Value = -1852.7144
$ws.Range("N74").Value = -5319.3333
$ws.Range("H77").Value = 2875.7646
$ws.Range("I77").Value = 2726.7144
$ws.Range("J77").Value = 3571.3333
$ws.Range("K77").Value = 13633.572
$ws.Range("L77").Value = 17856.6665
$ws.Range("M77").Value = -9265.572
$ws.Range("N77").Value = -26592.6665
$ws.Range("H116").Value = 1779
$ws.Range("I116").Value = 1779
$ws.Range("K116").Value = 1779
$ws.Range("M116").Value = 515
$ws.Range("H132").Value = 4846.242
$ws.Range("I132").Value = 4252.923
$ws.Range("J132").Value = 5231.9
$ws.Range("K132").Value = 12758.769
$ws.Range("L132").Value = 15695.7
$ws.Range("M132").Value = -10228.769
$ws.Range("N132").Value = -20755.7
$ws.Range("H136").Value = 2782.5151
$ws.Range("I136").Value = 2384.8096
$ws.Range("J136").Value = 3478.5
$ws.Range("K136").Value = 7154.4288
$ws.Range("L136").Value = 10435.5
$ws.Range("M136").Value = -4604.4288
$ws.Range("N136").Value = -15535.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1779
$ws.Range("I3").Value = 1779
$ws.Range("K3").Value = 1779
$ws.Range("M3").Value = -1665
$ws.Range("H105").Value = 25003234
$ws.Range("I105").Value = 31253142
$ws.Range("K105").Value = 31253142
$ws.Range("M105").Value = -31251395
$ws.Range("H134").Value = 3006.0356
$ws.Range("I134").Value = 3061.1428
$ws.Range("J134").Value = 2950.9285
$ws.Range("K134").Value = 9183.428400000001
$ws.Range("L134").Value = 8852.7855
$ws.Range("M134").Value = -6648.428400000001
$ws.Range("N134").Value = -13922.7855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 960.2
$ws.Range("K16").Value = 960.2
$ws.Range("M16").Value = -673.2
$ws.Range("H31").Value = 5330.9453
$ws.Range("I31").Value = 1294.2174
$ws.Range("K31").Value = 1294.2174
$ws.Range("M31").Value = -999.2174
$ws.Range("H34").Value = 5330.9453
$ws.Range("I34").Value = 1294.2174
$ws.Range("K34").Value = 1294.2174
$ws.Range("M34").Value = -1092.2174
$ws.Range("H58").Value = 1697.6364
$ws.Range("I58").Value = 2162
$ws.Range("J58").Value = 1561.0588
$ws.Range("K58").Value = 2162
$ws.Range("L58").Value = 1561.0588
$ws.Range("M58").Value = -1959
$ws.Range("N58").Value = -1967.0588
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 960.2
$ws.Range("K113").Value = 960.2
$ws.Range("M113").Value = 1209.8
$ws.Range("H122").Value = 1719.3914
$ws.Range("I122").Value = 1078.2858
$ws.Range("J122").Value = 1999.875
$ws.Range("K122").Value = 3234.8574
$ws.Range("L122").Value = 5999.625
$ws.Range("M122").Value = -784.8574000000003
$ws.Range("N122").Value = -10899.625
$ws.Range("H132").Value = 11907930
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 13891585
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 41674755
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -41679815
$ws.Range("H134").Value = 1798.0385
$ws.Range("I134").Value = 1634.1578
$ws.Range("K134").Value = 4902.4734
$ws.Range("M134").Value = -2367.4734
$ws.Range("H136").Value = 1697.6364
$ws.Range("I136").Value = 2162
$ws.Range("J136").Value = 1561.0588
$ws.Range("K136").Value = 6486
$ws.Range("L136").Value = 4683.1764
$ws.Range("M136").Value = -3936
$ws.Range("N136").Value = -9783.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 362.05884
$ws.Range("I5").Value = 362.05884
$ws.Range("K5").Value = 1086.17652
$ws.Range("M5").Value = -974.17652
$ws.Range("H14").Value = 220000060
$ws.Range("I14").Value = 220000060
$ws.Range("K14").Value = 660000180
$ws.Range("M14").Value = -660000007
$ws.Range("H98").Value = 271.125
$ws.Range("I98").Value = 250
$ws.Range("J98").Value = 283.8
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 851.4000000000001
$ws.Range("M98").Value = 748
$ws.Range("N98").Value = -3847.4
$ws.Range("H122").Value = 5782.2104
$ws.Range("I122").Value = 414.9091
$ws.Range("K122").Value = 3734.1819
$ws.Range("M122").Value = -1284.1819
$ws.Range("H135").Value = 362.05884
$ws.Range("I135").Value = 362.05884
$ws.Range("K135").Value = 3258.52956
$ws.Range("M135").Value = -723.5295599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2806.35
$ws.Range("I132").Value = 2346.4546
$ws.Range("K132").Value = 7039.3638
$ws.Range("M132").Value = -4509.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 12484.889
$ws.Range("I93").Value = 18160.834
$ws.Range("K93").Value = 18160.834
$ws.Range("M93").Value = -16912.834
$ws.Range("H132").Value = 3431
$ws.Range("I132").Value = 2192.3635
$ws.Range("J132").Value = 4566.4165
$ws.Range("K132").Value = 6577.0905
$ws.Range("L132").Value = 13699.2495
$ws.Range("M132").Value = -4047.0905
$ws.Range("N132").Value = -18759.2495
$ws.Range("H136").Value = 18521298
$ws.Range("J136").Value = 41669668
$ws.Range("L136").Value = 125009004
$ws.Range("N136").Value = -125014104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 487.22223
$ws.Range("I107").Value = 466.66666
$ws.Range("J107").Value = 497.5
$ws.Range("K107").Value = 1399.99998
$ws.Range("L107").Value = 1492.5
$ws.Range("M107").Value = 520.0000199999999
$ws.Range("N107").Value = -5332.5
$ws.Range("H132").Value = 4904640
$ws.Range("I132").Value = 3667.75
$ws.Range("J132").Value = 7577898
$ws.Range("K132").Value = 11003.25
$ws.Range("L132").Value = 22733694
$ws.Range("M132").Value = -8473.25
$ws.Range("N132").Value = -22738754
$ws.Range("H136").Value = 2670.3823
$ws.Range("I136").Value = 2525.818
$ws.Range("K136").Value = 7577.454000000001
$ws.Range("M136").Value = -5027.454000000001
